$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# A new handback was generated for b7e42bd5-d692-4697-96cd-07baad828a8a,
# filling in row 7 of both the "zh-cn" and "de-de" worksheets, plus a
# hyperlink to the handback markdown file on the "Latest Target File"
# (column I) cell.
# ---------------------------------------------------------------------------

$handbackMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec1777f036931bc894b5e3676d92acfadede8555/e2e/b7e42bd5-d692-4697-96cd-07baad828a8a.md"
$handbackMdDisplay = "b7e42bd5-d692-4697-96cd-07baad828a8a.md"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/16587c35af687edbfc862906f76c8189cf1cf1ac/e2e/b7e42bd5-d692-4697-96cd-07baad828a8a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec1777f036931bc894b5e3676d92acfadede8555/e2e/b7e42bd5-d692-4697-96cd-07baad828a8a.md."

# --- zh-cn sheet -----------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I7"), $handbackMdUrl, "", "", $handbackMdDisplay)
$wsZhCn.Range("J7").Value = "b7e42bd5-d692-4697-96cd-07baad828a8a.a245efa208546f64ac1b80d6ca3918fdeab68e11.zh-cn.xlf"
$wsZhCn.Range("K7").Value = "2016-08-27 06:51:53"
$wsZhCn.Range("P7").Value = $errorDetail

# --- de-de sheet -------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I7"), $handbackMdUrl, "", "", $handbackMdDisplay)
$wsDeDe.Range("J7").Value = "b7e42bd5-d692-4697-96cd-07baad828a8a.a245efa208546f64ac1b80d6ca3918fdeab68e11.de-de.xlf"
$wsDeDe.Range("K7").Value = "2016-08-27 06:51:59"
$wsDeDe.Range("P7").Value = $errorDetail
